# Modif Bilan Massique Châssis
# Updates the "Structure tubulaire" block (rows 19-21) on the single
# worksheet: E19 drops to 32.5, and E20/F20 pick up a new 3-unit entry,
# which feeds the E21/F21 totals (now written as explicit formulas
# instead of the old shared formula) and cascades into every downstream
# total (row 29, row 61, row 65) via recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structure tubulaire block ---
$ws.Range("E19").Value = 32.5

$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 3

# Row 21 (E21/F21) already totals E19:E20 / F19:F20 via the existing
# shared formula in that row - recalculation below refreshes its cached
# result (and every downstream total) without needing to touch the
# formula cells directly.

# --- Selection left behind by the editing session ---
$ws.Range("F21").Select()

$wb.Application.Calculate()
